# Auto-generated PowerShell COM-interop script
# Applies numeric corrections to the Leve profit-tracking sheets
# (ALC, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 3702.535
$ws.Cells.Item(15, 9).Value = 3702.535
$ws.Cells.Item(15, 11).Value = 11107.605
$ws.Cells.Item(15, 13).Value = -10938.605
$ws.Cells.Item(140, 8).Value = 37650
$ws.Cells.Item(140, 10).Value = 37650
$ws.Cells.Item(140, 12).Value = 37650
$ws.Cells.Item(140, 14).Value = -48010
$ws.Cells.Item(141, 8).Value = 6901.7617
$ws.Cells.Item(141, 9).Value = 2762.4285
$ws.Cells.Item(141, 10).Value = 15180.429
$ws.Cells.Item(141, 11).Value = 8287.2855
$ws.Cells.Item(141, 12).Value = 45541.287
$ws.Cells.Item(141, 13).Value = -3107.2855
$ws.Cells.Item(141, 14).Value = -55901.287

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 4152.636
$ws.Cells.Item(132, 9).Value = 3197.2856
$ws.Cells.Item(132, 10).Value = 5824.5
$ws.Cells.Item(132, 11).Value = 9591.856800000001
$ws.Cells.Item(132, 12).Value = 17473.5
$ws.Cells.Item(132, 13).Value = -7061.856800000001
$ws.Cells.Item(132, 14).Value = -22533.5
$ws.Cells.Item(134, 8).Value = 5724.8667
$ws.Cells.Item(134, 9).Value = 6343.9
$ws.Cells.Item(134, 10).Value = 4486.8
$ws.Cells.Item(134, 11).Value = 19031.7
$ws.Cells.Item(134, 12).Value = 13460.4
$ws.Cells.Item(134, 13).Value = -16496.7
$ws.Cells.Item(134, 14).Value = -18530.4

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1121.4783
$ws.Cells.Item(5, 9).Value = 1065.2106
$ws.Cells.Item(5, 10).Value = 1388.75
$ws.Cells.Item(5, 11).Value = 3195.6318
$ws.Cells.Item(5, 12).Value = 4166.25
$ws.Cells.Item(5, 13).Value = -3083.6318
$ws.Cells.Item(5, 14).Value = -4390.25
$ws.Cells.Item(131, 8).Value = 1195.5555
$ws.Cells.Item(131, 9).Value = 1950
$ws.Cells.Item(131, 10).Value = 980
$ws.Cells.Item(131, 11).Value = 5850
$ws.Cells.Item(131, 12).Value = 2940
$ws.Cells.Item(131, 13).Value = -810
$ws.Cells.Item(131, 14).Value = -13020
$ws.Cells.Item(135, 8).Value = 1121.4783
$ws.Cells.Item(135, 9).Value = 1065.2106
$ws.Cells.Item(135, 10).Value = 1388.75
$ws.Cells.Item(135, 11).Value = 9586.895400000001
$ws.Cells.Item(135, 12).Value = 12498.75
$ws.Cells.Item(135, 13).Value = -7051.895400000001
$ws.Cells.Item(135, 14).Value = -17568.75

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 5562222
$ws.Cells.Item(20, 10).Value = 7499.125
$ws.Cells.Item(20, 12).Value = 7499.125
$ws.Cells.Item(20, 14).Value = -7989.125
$ws.Cells.Item(24, 8).Value = 3000
$ws.Cells.Item(24, 10).Value = 3000
$ws.Cells.Item(24, 12).Value = 3000
$ws.Cells.Item(24, 14).Value = -3346
$ws.Cells.Item(63, 8).Value = 30000
$ws.Cells.Item(63, 10).Value = 30000
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 14).Value = -31372
$ws.Cells.Item(66, 8).Value = 30000
$ws.Cells.Item(66, 10).Value = 30000
$ws.Cells.Item(66, 12).Value = 90000
$ws.Cells.Item(66, 14).Value = -96864
$ws.Cells.Item(80, 8).Value = 2817.7046
$ws.Cells.Item(80, 9).Value = 2862.5
$ws.Cells.Item(80, 10).Value = 2807.75
$ws.Cells.Item(80, 11).Value = 2862.5
$ws.Cells.Item(80, 12).Value = 2807.75
$ws.Cells.Item(80, 13).Value = -1864.5
$ws.Cells.Item(80, 14).Value = -4803.75
$ws.Cells.Item(83, 8).Value = 2817.7046
$ws.Cells.Item(83, 9).Value = 2862.5
$ws.Cells.Item(83, 10).Value = 2807.75
$ws.Cells.Item(83, 11).Value = 14312.5
$ws.Cells.Item(83, 12).Value = 14038.75
$ws.Cells.Item(83, 13).Value = -9320.5
$ws.Cells.Item(83, 14).Value = -24022.75
$ws.Cells.Item(132, 8).Value = 2569.6843
$ws.Cells.Item(132, 9).Value = 2131.8333
$ws.Cells.Item(132, 10).Value = 3320.2856
$ws.Cells.Item(132, 11).Value = 6395.499899999999
$ws.Cells.Item(132, 12).Value = 9960.856800000001
$ws.Cells.Item(132, 13).Value = -3865.499899999999
$ws.Cells.Item(132, 14).Value = -15020.8568

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1399.5333
$ws.Cells.Item(16, 9).Value = 861
$ws.Cells.Item(16, 10).Value = 4900
$ws.Cells.Item(16, 11).Value = 861
$ws.Cells.Item(16, 12).Value = 4900
$ws.Cells.Item(16, 13).Value = -691
$ws.Cells.Item(16, 14).Value = -5240
$ws.Cells.Item(46, 8).Value = 1455.4
$ws.Cells.Item(46, 10).Value = 1764.8572
$ws.Cells.Item(46, 12).Value = 1764.8572
$ws.Cells.Item(46, 14).Value = -2140.8572
$ws.Cells.Item(68, 8).Value = 14099350
$ws.Cells.Item(68, 9).Value = 26027296
$ws.Cells.Item(68, 10).Value = 2687.182
$ws.Cells.Item(68, 11).Value = 26027296
$ws.Cells.Item(68, 12).Value = 2687.182
$ws.Cells.Item(68, 13).Value = -26026547
$ws.Cells.Item(68, 14).Value = -4185.182
$ws.Cells.Item(71, 8).Value = 14099350
$ws.Cells.Item(71, 9).Value = 26027296
$ws.Cells.Item(71, 10).Value = 2687.182
$ws.Cells.Item(71, 11).Value = 130136480
$ws.Cells.Item(71, 12).Value = 13435.91
$ws.Cells.Item(71, 13).Value = -130132736
$ws.Cells.Item(71, 14).Value = -20923.91
$ws.Cells.Item(82, 8).Value = 3593.25
$ws.Cells.Item(82, 9).Value = 3600
$ws.Cells.Item(82, 10).Value = 3591
$ws.Cells.Item(82, 11).Value = 3600
$ws.Cells.Item(82, 12).Value = 3591
$ws.Cells.Item(82, 13).Value = -3239
$ws.Cells.Item(82, 14).Value = -4313
$ws.Cells.Item(85, 8).Value = 3593.25
$ws.Cells.Item(85, 9).Value = 3600
$ws.Cells.Item(85, 10).Value = 3591
$ws.Cells.Item(85, 11).Value = 3600
$ws.Cells.Item(85, 12).Value = 3591
$ws.Cells.Item(85, 13).Value = -2352
$ws.Cells.Item(85, 14).Value = -6087
$ws.Cells.Item(93, 8).Value = 1797.2727
$ws.Cells.Item(93, 9).Value = 1668.5714
$ws.Cells.Item(93, 11).Value = 1668.5714
$ws.Cells.Item(93, 13).Value = -420.5714
$ws.Cells.Item(132, 8).Value = 3693364.5
$ws.Cells.Item(132, 9).Value = 6970266
$ws.Cells.Item(132, 10).Value = 6850.375
$ws.Cells.Item(132, 11).Value = 20910798
$ws.Cells.Item(132, 12).Value = 20551.125
$ws.Cells.Item(132, 13).Value = -20908268
$ws.Cells.Item(132, 14).Value = -25611.125

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 17333.334
$ws.Cells.Item(54, 10).Value = 17333.334
$ws.Cells.Item(54, 12).Value = 17333.334
$ws.Cells.Item(54, 14).Value = -18373.334
$ws.Cells.Item(62, 8).Value = 5456.9165
$ws.Cells.Item(62, 9).Value = 3660.8
$ws.Cells.Item(62, 10).Value = 6739.857
$ws.Cells.Item(62, 11).Value = 3660.8
$ws.Cells.Item(62, 12).Value = 6739.857
$ws.Cells.Item(62, 13).Value = -3036.8
$ws.Cells.Item(62, 14).Value = -7987.857
$ws.Cells.Item(65, 8).Value = 5456.9165
$ws.Cells.Item(65, 9).Value = 3660.8
$ws.Cells.Item(65, 10).Value = 6739.857
$ws.Cells.Item(65, 11).Value = 18304
$ws.Cells.Item(65, 12).Value = 33699.285
$ws.Cells.Item(65, 13).Value = -15184
$ws.Cells.Item(65, 14).Value = -39939.285
$ws.Cells.Item(107, 8).Value = 1063.6
$ws.Cells.Item(107, 9).Value = 1002
$ws.Cells.Item(107, 10).Value = 1079
$ws.Cells.Item(107, 11).Value = 3006
$ws.Cells.Item(107, 12).Value = 3237
$ws.Cells.Item(107, 13).Value = -1086
$ws.Cells.Item(107, 14).Value = -7077
